# TetraPlex deck 0.54 - dark mode
#
# Fix a typo in the title of the "Tetraplex Kognitiv... : AI Services"
# slide (slide 11, title placeholder / shape 1):
#
#   "Tetraplex Kognitivdiemste:  AI Services"
#              ^^^^^^^^^^^^^^^^
#        -> "Tetraplex Kognitivdienste:  AI Services"
#
# Only the misspelled word is touched; the rest of the run is left
# untouched (PowerPoint will naturally split the original single run
# into "Tetraplex " / "Kognitivdienste" / ":  AI Services" around the
# edited substring).

$p = $ppt.ActivePresentation

$needle = "Kognitivdiemste"
$fixed = "Kognitivdienste"

$target = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -match [regex]::Escape($needle)) {
            $target = $shp.TextFrame.TextRange
            break
        }
    }
    if ($target -ne $null) { break }
}

# Fall back to the known location (slide 11, shape 1) if the search
# above didn't find anything, e.g. because the text was already fixed.
if ($target -eq $null) {
    $target = $p.Slides.Item(11).Shapes.Item(1).TextFrame.TextRange
}

$fullText = $target.Text
$idx = $fullText.IndexOf($needle)
if ($idx -ge 0) {
    # Characters() is 1-based.
    $sub = $target.Characters($idx + 1, $needle.Length)
    $sub.Text = $fixed
}
